# Rename the three header/footer logo pictures:
#   - Footer "first page" (footer1.xml, Pearson logo, docPr id="3")  image1.png -> image2.png
#   - Footer "default"    (footer2.xml, Pearson logo, docPr id="2")  image1.png -> image2.png
#   - Header "first page" (header1.xml, BTec logo,    docPr id="1")  image2.jpg -> image1.jpg
#
# InlineShape has no settable .Name in Word's object model, so each
# picture is temporarily promoted to a floating Shape (which does
# expose .Name), renamed, then converted back to an InlineShape so the
# surrounding <w:drawing><wp:inline> markup/run layout is preserved.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-LogoPicture($range, $newName) {
    $inlineShape = $range.InlineShapes.Item(1)
    $shape = $inlineShape.ConvertToShape()
    $shape.Name = $newName
    [void]$shape.ConvertToInlineShape()
}

# Footer, first page variant (footer1.xml) - Pearson logo id="3"
$footerFirst = $sec.Footers.Item(2)
Rename-LogoPicture $footerFirst.Range "image2.png"

# Footer, default variant (footer2.xml) - Pearson logo id="2"
$footerDefault = $sec.Footers.Item(1)
Rename-LogoPicture $footerDefault.Range "image2.png"

# Header, first page variant (header1.xml) - BTec logo id="1"
$headerFirst = $sec.Headers.Item(2)
Rename-LogoPicture $headerFirst.Range "image1.jpg"
